$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width tweaks (auto-fit-like adjustments from new header/content widths) ---
$ws.Columns.Item(3).ColumnWidth = 13.333333333333334
$ws.Columns.Item(4).ColumnWidth = 14.333333333333334
$ws.Columns.Item(7).ColumnWidth = 15.333333333333334
$ws.Columns.Item(8).ColumnWidth = 14.833333333333334
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
$ws.Columns.Item(11).ColumnWidth = 15.0

# --- Roll the nowcast window forward by one quarter (new row dates + values) ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-09-30"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = 0.27217399496364714
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025-10-15"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = 0.2711959182004798
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = -0.005674509449603661
$ws.Range("E3").Value = -0.000224570284440243
$ws.Range("F3").Value = -0.0005506061464371004
$ws.Range("G3").Value = 0.0001392507223496883
$ws.Range("H3").Value = 0.000026821926490800298
$ws.Range("I3").Value = 0.000055597310817432816
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0.00007642868917062673

$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2025-10-30"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = 0.33171928741532447
$ws.Range("C4").Value = 0.019820956160900682
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = -0.000003976380617616612
$ws.Range("F4").Value = 0.000002574330223397678
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.00013507365155167425
$ws.Range("I4").Value = -0.00039358390289919354
$ws.Range("J4").Value = 0.0029182491864013726
$ws.Range("K4").Value = -0.00028317984103581884

$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2025-11-15"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = 0.3285216845327023
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = -0.002325679714222684
$ws.Range("E5").Value = -0.0009719723648560505
$ws.Range("F5").Value = 0.0015852282811317589
$ws.Range("G5").Value = -0.000494387337057782
$ws.Range("H5").Value = 0.00013724234887271678
$ws.Range("I5").Value = -0.00010901462737906834
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = -0.001831087083280114

$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2025-11-30"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = 0.18158653251933543
$ws.Range("C6").Value = -0.014758525507606514
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0.000024074485978141254
$ws.Range("F6").Value = 0.00023079686581279588
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0.00013652148101326596
$ws.Range("I6").Value = -0.006705332678328006
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = -0.00004318556087806735

$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2025-12-15"
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = 0.1816367577376321
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = -0.005518875670399108
$ws.Range("E7").Value = -0.0006396056111778178
$ws.Range("F7").Value = 0.002599780350363929
$ws.Range("G7").Value = 0.0031114559852409235
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = -0.0007480262824329509
